# Insert a new weekly price record as row 48, pushing all subsequent
# rows (old 48..154) down by one (new 49..155).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("48").Insert()

$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44497
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 100112030
$ws.Range("G48").Value = "Poroto granado"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 55
$ws.Range("K48").Value = 40000
$ws.Range("L48").Value = 45000
$ws.Range("M48").Value = 42273
$ws.Range("N48").Value = "`$/malla 25 kilos"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value = 1691
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
